$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 10
$ws.Range("F10").Value = 564
$ws.Range("H10").Value = 659

# Row 11
$ws.Range("F11").Value = 380
$ws.Range("H11").Value = 445

# Row 12
$ws.Range("F12").Value = 615
$ws.Range("H12").Value = 701

# Row 15
$ws.Range("G15").Value = 52
$ws.Range("H15").Value = 182

# Row 23
$ws.Range("F23").Value = 154
$ws.Range("H23").Value = 206

# Row 24
$ws.Range("F24").Value = 224
$ws.Range("H24").Value = 254

# Row 26
$ws.Range("F26").Value = 185
$ws.Range("H26").Value = 210

# Row 28
$ws.Range("F28").Value = 157
$ws.Range("H28").Value = 209

# Row 41
$ws.Range("F41").Value = 319
$ws.Range("H41").Value = 411

# Row 42
$ws.Range("F42").Value = 376
$ws.Range("H42").Value = 437

# Row 47
$ws.Range("F47").Value = 433
$ws.Range("H47").Value = 525

# Row 48
$ws.Range("F48").Value = 186
$ws.Range("H48").Value = 230

# Row 51
$ws.Range("F51").Value = 159
$ws.Range("H51").Value = 233

$wb.Save()
